$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 133 (shifts existing rows 133:163 down to 134:164)
$ws.Rows.Item(133).Insert()

# Populate the newly inserted row 133 with the new record
$ws.Cells.Item(133, 1).Value = 3
$ws.Cells.Item(133, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(133, 3).Value = "Coquimbo"
$ws.Cells.Item(133, 4).Value = 44543
$ws.Cells.Item(133, 5).Value = 5
$ws.Cells.Item(133, 6).Value = 100112010
$ws.Cells.Item(133, 7).Value = "Achicoria"
$ws.Cells.Item(133, 8).Value = "Sin especificar"
$ws.Cells.Item(133, 9).Value = "Primera"
$ws.Cells.Item(133, 10).Value = 50
$ws.Cells.Item(133, 11).Value = 5500
$ws.Cells.Item(133, 12).Value = 5500
$ws.Cells.Item(133, 13).Value = 5500
$ws.Cells.Item(133, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(133, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(133, 16).Value = 344
$ws.Cells.Item(133, 17).Value = 16
$ws.Cells.Item(133, 18).Value = "Hortaliza"
